# "Updated Stats - Wk 13"
# Refresh the weekly fantasy-football stats sheet: Z/AN columns (week's
# matchup score, previously blank/"NA") now have real numbers for the
# week that just completed, a handful of season-to-date counters
# (Acquisitions/Drops/Wins/Losses) ticked up, and the WinPct column's
# display format gained a third decimal place. H/I/M are formulas so
# they recalculate on their own once the inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("J2").Value = 7
$ws.Range("Z2").Value = 166
$ws.Range("AN2").Value = 87

# --- Row 3 ---
$ws.Range("F3").Value = 57
$ws.Range("K3").Value = 8
$ws.Range("Z3").Value = 87
$ws.Range("AN3").Value = 166

# --- Row 4 ---
$ws.Range("F4").Value = 55
$ws.Range("J4").Value = 8
$ws.Range("Z4").Value = 182.5
$ws.Range("AN4").Value = 136

# --- Row 5 ---
$ws.Range("D5").Value = 42
$ws.Range("E5").Value = 42
$ws.Range("F5").Value = 89
$ws.Range("K5").Value = 7
$ws.Range("Z5").Value = 117.5
$ws.Range("AN5").Value = 137

# --- Row 6 ---
$ws.Range("F6").Value = 71
$ws.Range("K6").Value = 8
$ws.Range("Z6").Value = 61.5
$ws.Range("AN6").Value = 134

# --- Row 7 ---
$ws.Range("F7").Value = 55
$ws.Range("J7").Value = 5
$ws.Range("Z7").Value = 134
$ws.Range("AN7").Value = 61.5

# --- Row 8 ---
$ws.Range("J8").Value = 8
$ws.Range("Z8").Value = 142
$ws.Range("AN8").Value = 69.5

# --- Row 9 ---
$ws.Range("K9").Value = 8
$ws.Range("Z9").Value = 69.5
$ws.Range("AN9").Value = 142

# --- Row 10 ---
$ws.Range("F10").Value = 41
$ws.Range("K10").Value = 5
$ws.Range("Z10").Value = 136
$ws.Range("AN10").Value = 182.5

# --- Row 11 ---
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 31
$ws.Range("F11").Value = 70
$ws.Range("J11").Value = 7
$ws.Range("Z11").Value = 137
$ws.Range("AN11").Value = 117.5

# WinPct column now shows three decimal places instead of two.
$ws.Range("M2:M11").NumberFormat = "0.000"

# Stale sort-state bookmark (below the real data) shifted up two rows.
$ws.Sort.SortFields.Clear()
$ws.Sort.SetRange($ws.Range("B22:G31"))
$null = $ws.Sort.SortFields.Add($ws.Range("B21"))
$ws.Sort.Apply()

# Leave the cursor where the editor last left it.
$ws.Range("AN12").Select()
